$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet named "Sheet2" positioned right after Sheet1
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

# Enter the dynamic-array TRANSPOSE formula that spills from A1 into A1:AE3,
# reproducing (transposing) Sheet1's A1:C31 range (31 rows x 3 cols -> 3 rows x 31 cols)
$ws2.Range("A1").Formula2 = "=TRANSPOSE(Sheet1!A1:C31)"

# Match the approximate column sizing used in the authored sheet
$ws2.Columns.Item(1).ColumnWidth = 37.166666666666664
$ws2.Columns.Item(3).ColumnWidth = 14.858072916666666

# Match the recorded selection / cursor position on Sheet2
$ws2.Range("A11").Select() | Out-Null

# Keep Sheet1 as the active/visible tab, as in the target workbook
$ws1.Activate() | Out-Null
